# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 2.598097515653722)
    3 = @(1.445647641019636, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 16.33547987844073)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    5 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 13.86384647080068, 17.08608867836142)
    6 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 17.65757632934944)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
